$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet) - F column "想去人数" (want-to-go count) updates
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F4").Value = 106
$wsExhibit.Range("F5").Value = 141
$wsExhibit.Range("F6").Value = 1359
$wsExhibit.Range("F7").Value = 1594
$wsExhibit.Range("F8").Value = 349
$wsExhibit.Range("F9").Value = 459
$wsExhibit.Range("F11").Value = 188
$wsExhibit.Range("F16").Value = 326
$wsExhibit.Range("F18").Value = 1778
$wsExhibit.Range("F22").Value = 703
$wsExhibit.Range("F24").Value = 349
$wsExhibit.Range("F25").Value = 4307
$wsExhibit.Range("F26").Value = 15
$wsExhibit.Range("F28").Value = 1139
$wsExhibit.Range("F29").Value = 498
$wsExhibit.Range("F31").Value = 668
$wsExhibit.Range("F33").Value = 331
$wsExhibit.Range("F35").Value = 172

# Sheet "全部类型" (4th sheet) - F column "想去人数" (want-to-go count) updates
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F4").Value = 106
$wsAll.Range("F5").Value = 141
$wsAll.Range("F6").Value = 1359
$wsAll.Range("F7").Value = 1594
$wsAll.Range("F8").Value = 349
$wsAll.Range("F9").Value = 460
$wsAll.Range("F11").Value = 188
$wsAll.Range("F16").Value = 326
$wsAll.Range("F18").Value = 1778
$wsAll.Range("F22").Value = 703
$wsAll.Range("F24").Value = 349
$wsAll.Range("F25").Value = 4307
$wsAll.Range("F26").Value = 15
$wsAll.Range("F28").Value = 1139
$wsAll.Range("F29").Value = 498
$wsAll.Range("F31").Value = 668
$wsAll.Range("F33").Value = 331
$wsAll.Range("F35").Value = 172
